$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new row (row 8) that mirrors row 7 ("Quy trình đào tạo"),
#     but with its own rich-text content coming from the text editor ---

# A8: plain numeric id
$ws.Range("A8").Value = 8

# B8: project name (reuses the existing "Quy trình đào tạo" string)
$ws.Range("B8").Value = "Quy trình đào tạo"

# C8 / F8 need to stay plain TEXT (like row 7's C7/F7), not get auto-converted
# into a date serial / number by the COM layer. Force text via NumberFormat,
# assign, then drop the number-format override again so the cell ends up
# with the default style (matches how C7/F7 look - no explicit style).
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "09/09/2022"
$ws.Range("C8").ClearFormats()

# D8: rich HTML coming from the WYSIWYG text editor
$html = "<p><strong style=`"font-size: 18px;`">Đây là quy trình đào tạo</strong></p><table style=`"border-collapse:collapse;width: 100%;`"><tbody>`r`n<tr>`r`n`t<td style=`"width: 50%;`">Hello</td>`r`n`t<td style=`"width: 50%;`">Test</td></tr>`r`n<tr>`r`n`t<td style=`"width: 50%;`"><br></td>`r`n`t<td style=`"width: 50%;`"><br></td></tr></tbody></table>"
$ws.Range("D8").Value = $html

# The embedded line breaks make the engine auto-grow the row; put it back to
# the sheet's normal (non-custom) height right away so it doesn't linger and
# corrupt style bookkeeping for cells touched afterwards.
$ws.Rows.Item(8).AutoFit()

# E8: status (reuses existing "Chưa duyệt" string)
$ws.Range("E8").Value = "Chưa duyệt"

# F8: kept as text "123" (like F7's text "1233"), same trick as C8
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "123"
$ws.Range("F8").ClearFormats()

# Touch the header/footer so the sheet gains a (blank) <headerFooter/> entry,
# matching the resave performed by the originating tool.
$ws.PageSetup.CenterHeader = ""
